# Android-code.xlsx update:
# - add three new rows (Web View / startActivity guard / AlertDialog) to "Kotlin" sheet
# - add two new rows (multi-line EditText / AlertDialog button color style) to "UI" sheet
# - make "UI" the active/selected sheet, scrolled so row 7 is visible with B7 selected

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Kotlin" (sheet1.xml): append rows 5-7
# ---------------------------------------------------------------------------
$kotlin = $wb.Worksheets.Item("Kotlin")

$webViewCode = @"
webView.settings.javaScriptEnabled = true
        webView.settings.builtInZoomControls = true
        webView.settings.textZoom = 100
        webView.webViewClient = WebViewClient()
        webView.loadUrl(url)
"@ -replace "`n", "`r`n"

$kotlin.Range("A5").Value = $webViewCode
$kotlin.Range("A5").HorizontalAlignment = -4131
$kotlin.Range("A5").VerticalAlignment = -4108
$kotlin.Range("A5").WrapText = $true

$kotlin.Range("B5").Value = "Web View"
$kotlin.Range("B5").HorizontalAlignment = -4108
$kotlin.Range("B5").VerticalAlignment = -4108
$kotlin.Range("B5").WrapText = $false

$kotlin.Rows.Item(5).RowHeight = 75

$startActivityCode = @"
if (intent.resolveActivity(context.packageManager) != null)
                        context.startActivity(intent)
"@ -replace "`n", "`r`n"

$kotlin.Range("A6").Value = $startActivityCode
$kotlin.Range("A6").HorizontalAlignment = -4131
$kotlin.Range("A6").VerticalAlignment = -4108
$kotlin.Range("A6").WrapText = $true

$startActivityTitle = @"
اگر اکتیویتیش وجود نداشته باشه استارتش نمیکنه
از try/cach هم میتونیم استفاده کنیم
"@ -replace "`n", "`r`n"

$kotlin.Range("B6").Value = $startActivityTitle
$kotlin.Range("B6").VerticalAlignment = -4108
$kotlin.Range("B6").WrapText = $true

$kotlin.Rows.Item(6).RowHeight = 30

$alertDialogCode = @"
val view = LayoutInflater.from(context).inflate(R.layout.dialog_call_me, null)
                AlertDialog.Builder(context)
                    .setView(view)
                    .show()
"@ -replace "`n", "`r`n"

$kotlin.Range("A7").Value = $alertDialogCode
$kotlin.Range("A7").HorizontalAlignment = -4131
$kotlin.Range("A7").VerticalAlignment = -4108
$kotlin.Range("A7").WrapText = $true

$kotlin.Range("B7").Value = "AlertDialog"
$kotlin.Range("B7").HorizontalAlignment = -4108
$kotlin.Range("B7").VerticalAlignment = -4108
$kotlin.Range("B7").WrapText = $false

$kotlin.Rows.Item(7).RowHeight = 90

# ---------------------------------------------------------------------------
# Sheet "UI" (sheet2.xml): append rows 5-6
# ---------------------------------------------------------------------------
$ui = $wb.Worksheets.Item("UI")

$multilineEditText = @"

    <androidx.appcompat.widget.AppCompatEditText
        android:id="@+id/edt_text_dialog_failure_report"
        android:layout_width="0dp"
        android:layout_height="wrap_content"
        android:ems="10"
        android:gravity="top"
        android:hint="@string/text"
        android:inputType="textMultiLine"
        android:lines="3"
        android:maxLines="10"
        android:scrollbars="vertical"/>
        
"@ -replace "`n", "`r`n"

$ui.Range("A5").Value = $multilineEditText
$ui.Range("B5").Value = "تکست چند خطی"
$ui.Rows.Item(5).RowHeight = 195

$dialogThemeStyle = @"
<style name="DialogTheme" parent="Theme.AppCompat.Light.Dialog.Alert">
        <item name="buttonBarNegativeButtonStyle">@style/NegativeButton</item>
        <item name="buttonBarPositiveButtonStyle">@style/PositiveButton</item>
    </style>
 <style name="PositiveButton" parent="Widget.AppCompat.Button.ButtonBar.AlertDialog">
        <item name="android:textColor">@color/colorBlack</item>
        <item name="android:textSize">16sp</item>
    </style>
<style name="NegativeButton" parent="Widget.AppCompat.Button.ButtonBar.AlertDialog">
        <item name="android:textColor">@color/colorBlack</item>
        <item name="android:textSize">16sp</item>
    </style>
"@ -replace "`n", "`r`n"

$ui.Range("A6").Value = $dialogThemeStyle
$ui.Range("B6").Value = "برای رنگ به متن Alert Dilog"
$ui.Rows.Item(6).RowHeight = 255

# ---------------------------------------------------------------------------
# Make "UI" the active sheet / tab, scrolled to show row 6 onward, B7 selected
# ---------------------------------------------------------------------------
$ui.Activate()
$ui.Range("B7").Select()
$excel.ActiveWindow.ScrollRow = 6
